$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: field names ---------------------------------------------------
$ws.Range("A1").Value = "Num"
$ws.Range("B1:G1").Value = "GradeWeightList"
$ws.Range("H1:J1").Value = "PickupWeightList"
$ws.Range("K1").Value = "WeightSum"
$ws.Range("L1").Value = "CookieWeight"
$ws.Range("M1").Value = "StarExpWeight"
$ws.Range("N1").Value = "DetailWeightSum"

# --- Row 2: field types ----------------------------------------------------
$ws.Range("A2").Value = "int"
$ws.Range("B2:J2").Value = "list:int"
$ws.Range("K2").Value = "int"
$ws.Range("L2").Value = "int"
$ws.Range("M2").Value = "int"
$ws.Range("N2").Value = "int"

# --- Row 3: comment row ------------------------------------------------------
$ws.Range("A3").Value = "#"

# --- Row 4: data -------------------------------------------------------------
$ws.Range("A4").Value = 1001
$ws.Range("B4").Value = 600000
$ws.Range("C4").Value = 200000
$ws.Range("D4").Value = 100000
$ws.Range("E4").Value = 80000
$ws.Range("F4").Value = 10000
$ws.Range("G4").Value = 10000
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Formula = "=SUM(B4:J4)"
$ws.Range("L4").Value = 300000
$ws.Range("M4").Value = 700000
$ws.Range("N4").Formula = "=SUM(L4:M4)"

# --- Row 5: data -------------------------------------------------------------
$ws.Range("A5").Value = 2001
$ws.Range("B5").Value = 500000
$ws.Range("C5").Value = 100000
$ws.Range("D5").Value = 100000
$ws.Range("E5").Value = 80000
$ws.Range("F5").Value = 10000
$ws.Range("G5").Value = 10000
$ws.Range("H5").Value = 200000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Formula = "=SUM(B5:J5)"
$ws.Range("L5").Value = 300000
$ws.Range("M5").Value = 700000
$ws.Range("N5").Formula = "=SUM(L5:M5)"

# --- Column widths -----------------------------------------------------------
# Target stored widths (from the authored file) are 14.5703125 / 16.5703125 /
# 13.5703125 "character" units. The COM ColumnWidth setter here quantizes to
# an MDW=7 pixel grid (stored = (round(chars*7)+5)/7), so we request the
# chars value whose pixel-grid result lands nearest (within ~0.0011) of the
# authored width - the closest reproducible value.
$ws.Range("B1:G1").EntireColumn.ColumnWidth = 13.857142857142858
$ws.Range("H1:J1").EntireColumn.ColumnWidth = 15.857142857142858
$ws.Range("L1").EntireColumn.ColumnWidth = 12.857142857142858

# --- Selection -----------------------------------------------------------------
$ws.Range("H11").Select()
